$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 1554.359380824244
$ws.Range("B3").Value = 2005.105169756474
$ws.Range("B4").Value = 2523.232787144853
$ws.Range("B5").Value = 2788.772826141984
$ws.Range("B6").Value = 3079.659038269087
$ws.Range("B7").Value = 3296.728314284335
$ws.Range("B8").Value = 3521.459651329038
$ws.Range("B9").Value = 3709.498940577651
$ws.Range("B10").Value = 3933.15926532961
$ws.Range("B11").Value = 4027.577236152459
$ws.Range("B12").Value = 4194.453486219991
$ws.Range("B13").Value = 4264.331556590987
$ws.Range("B14").Value = 4374.067528604009
$ws.Range("B15").Value = 4487.385843203305
$ws.Range("B16").Value = 4576.628211015293
$ws.Range("B17").Value = 4626.525892366498
$ws.Range("B18").Value = 4665.233919472924
$ws.Range("B19").Value = 4695.277210699368
$ws.Range("B20").Value = 4743.863204282354
$ws.Range("B21").Value = 4755.222104715553
$ws.Range("B22").Value = 4773.863366161764
$ws.Range("B23").Value = 4764.742751547747
$ws.Range("B24").Value = 4779.311739273581
$ws.Range("B25").Value = 4755.846064415698
$ws.Range("B26").Value = 4708.867678279513
$ws.Range("B27").Value = 4682.215453613936
$ws.Range("B28").Value = 4652.618428692093
$ws.Range("B29").Value = 4599.074200781597
$ws.Range("B30").Value = 4537.17695691405
$ws.Range("B31").Value = 4488.640088597283
$ws.Range("B32").Value = 4420.682544029601
$ws.Range("B33").Value = 4391.841064129753
$ws.Range("B34").Value = 4300.46656922658
$ws.Range("B35").Value = 4187.720798297653
$ws.Range("B36").Value = 4104.633814792059
$ws.Range("B37").Value = 4008.160257769249
$ws.Range("B38").Value = 3835.890485365039
$ws.Range("B39").Value = 3775.843548412584
$ws.Range("B40").Value = 3652.012726918672
$ws.Range("B41").Value = 3541.271316119653
$ws.Range("B42").Value = 3465.212594074099
$ws.Range("B43").Value = 3357.511238168855
$ws.Range("B44").Value = 3208.694170771299
$ws.Range("B45").Value = 3122.890725073222
$ws.Range("B46").Value = 3020.4286909934
$ws.Range("B47").Value = 2907.454447358568
$ws.Range("B48").Value = 2722.580345836597
$ws.Range("B49").Value = 2643.834802528551
$ws.Range("B50").Value = 2490.081011636917
$ws.Range("B51").Value = 2279.523451392301
$ws.Range("B52").Value = 1919.020804917689
$ws.Range("B53").Value = 1644.7244712366
$ws.Range("B54").Value = 1575.406007724059
$ws.Range("B55").Value = 1546.550633036336
$ws.Range("B56").Value = 1277.343174404287
$ws.Range("B57").Value = 1116.674075305931
$ws.Range("B58").Value = 1002.166736952312
$ws.Range("B59").Value = 973.3299723037687
$ws.Range("B60").Value = 964.452203373652
$ws.Range("B61").Value = 960.9871111218379
$ws.Range("B62").Value = 957.1438016332409
